$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29379.75
$ws.Range("I21").Value = 25005.428
$ws.Range("J21").Value = 60000
$ws.Range("K21").Value = 25005.428
$ws.Range("L21").Value = 60000
$ws.Range("M21").Value = -24537.428
$ws.Range("N21").Value = -60936
$ws.Range("H23").Value = 29379.75
$ws.Range("I23").Value = 25005.428
$ws.Range("J23").Value = 60000
$ws.Range("K23").Value = 25005.428
$ws.Range("L23").Value = 60000
$ws.Range("M23").Value = -24771.428
$ws.Range("N23").Value = -60468
$ws.Range("H135").Value = 55557670
$ws.Range("I135").Value = 2379.125
$ws.Range("K135").Value = 21412.125
$ws.Range("M135").Value = -18877.125
$ws.Range("H137").Value = 3478.2712
$ws.Range("I137").Value = 895.5
$ws.Range("J137").Value = 4005.3674
$ws.Range("K137").Value = 2686.5
$ws.Range("L137").Value = 12016.1022
$ws.Range("M137").Value = -136.5
$ws.Range("N137").Value = -17116.1022

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1050.5
$ws.Range("J4").Value = 1001
$ws.Range("L4").Value = 1001
$ws.Range("N4").Value = -1233
$ws.Range("H10").Value = 100003.336
$ws.Range("J10").Value = 100003.336
$ws.Range("L10").Value = 100003.336
$ws.Range("N10").Value = -100343.336
$ws.Range("H23").Value = 80006
$ws.Range("I23").Value = 80006
$ws.Range("K23").Value = 80006
$ws.Range("M23").Value = -79747
$ws.Range("H37").Value = 37496.668
$ws.Range("J37").Value = 37496.668
$ws.Range("L37").Value = 37496.668
$ws.Range("N37").Value = -38042.668
$ws.Range("H44").Value = 36278.082
$ws.Range("J44").Value = 36278.082
$ws.Range("L44").Value = 36278.082
$ws.Range("N44").Value = -37254.082
$ws.Range("H55").Value = 49222.5
$ws.Range("J55").Value = 49222.5
$ws.Range("L55").Value = 49222.5
$ws.Range("N55").Value = -49852.5
$ws.Range("H61").Value = 2794.64
$ws.Range("I61").Value = 1363.3914
$ws.Range("J61").Value = 4013.8518
$ws.Range("K61").Value = 1363.3914
$ws.Range("L61").Value = 4013.8518
$ws.Range("M61").Value = -1151.3914
$ws.Range("N61").Value = -4437.8518
$ws.Range("H74").Value = 1315
$ws.Range("I74").Value = 818
$ws.Range("J74").Value = 2159.9
$ws.Range("K74").Value = 818
$ws.Range("L74").Value = 2159.9
$ws.Range("M74").Value = 56
$ws.Range("N74").Value = -3907.9
$ws.Range("H77").Value = 1315
$ws.Range("I77").Value = 818
$ws.Range("J77").Value = 2159.9
$ws.Range("K77").Value = 4090
$ws.Range("L77").Value = 10799.5
$ws.Range("M77").Value = 278
$ws.Range("N77").Value = -19535.5
$ws.Range("H80").Value = 54095
$ws.Range("J80").Value = 54095
$ws.Range("L80").Value = 54095
$ws.Range("N80").Value = -56091
$ws.Range("H83").Value = 54095
$ws.Range("J83").Value = 54095
$ws.Range("L83").Value = 162285
$ws.Range("N83").Value = -172269
$ws.Range("H102").Value = 35100
$ws.Range("I102").Value = 1400
$ws.Range("J102").Value = 51950
$ws.Range("K102").Value = 1400
$ws.Range("L102").Value = 51950
$ws.Range("M102").Value = 222
$ws.Range("N102").Value = -55194
$ws.Range("H109").Value = 46518
$ws.Range("J109").Value = 46518
$ws.Range("L109").Value = 46518
$ws.Range("N109").Value = -49292
$ws.Range("H122").Value = 1930.6666
$ws.Range("I122").Value = 1974.6666
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 5923.9998
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -3473.9998
$ws.Range("N122").Value = -9899.9998
$ws.Range("H136").Value = 2794.64
$ws.Range("I136").Value = 1363.3914
$ws.Range("J136").Value = 4013.8518
$ws.Range("K136").Value = 4090.1742
$ws.Range("L136").Value = 12041.5554
$ws.Range("M136").Value = -1540.1742
$ws.Range("N136").Value = -17141.5554

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 33981.668
$ws.Range("J19").Value = 33981.668
$ws.Range("L19").Value = 33981.668
$ws.Range("N19").Value = -34327.668
$ws.Range("H108").Value = 37283.668
$ws.Range("J108").Value = 37283.668
$ws.Range("L108").Value = 37283.668
$ws.Range("N108").Value = -44963.668
$ws.Range("H134").Value = 4276.711
$ws.Range("I134").Value = 2810.1428
$ws.Range("J134").Value = 4773.4517
$ws.Range("K134").Value = 8430.428400000001
$ws.Range("L134").Value = 14320.3551
$ws.Range("M134").Value = -5895.428400000001
$ws.Range("N134").Value = -19390.3551

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3992.674
$ws.Range("I31").Value = 1306.4546
$ws.Range("K31").Value = 1306.4546
$ws.Range("M31").Value = -1011.4546
$ws.Range("H34").Value = 3992.674
$ws.Range("I34").Value = 1306.4546
$ws.Range("K34").Value = 1306.4546
$ws.Range("M34").Value = -1104.4546
$ws.Range("H58").Value = 4538.8057
$ws.Range("I58").Value = 4989.1724
$ws.Range("J58").Value = 2673
$ws.Range("K58").Value = 4989.1724
$ws.Range("L58").Value = 2673
$ws.Range("M58").Value = -4786.1724
$ws.Range("N58").Value = -3079
$ws.Range("H103").Value = 15203.2
$ws.Range("I103").Value = 5629.6
$ws.Range("J103").Value = 19990
$ws.Range("K103").Value = 5629.6
$ws.Range("L103").Value = 19990
$ws.Range("M103").Value = -4457.6
$ws.Range("N103").Value = -22334
$ws.Range("H136").Value = 4538.8057
$ws.Range("I136").Value = 4989.1724
$ws.Range("J136").Value = 2673
$ws.Range("K136").Value = 14967.5172
$ws.Range("L136").Value = 8019
$ws.Range("M136").Value = -12417.5172
$ws.Range("N136").Value = -13119

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3757.2188
$ws.Range("I113").Value = 4951.478
$ws.Range("J113").Value = 705.2222
$ws.Range("K113").Value = 14854.434
$ws.Range("L113").Value = 2115.6666
$ws.Range("M113").Value = -12684.434
$ws.Range("N113").Value = -6455.6666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 50014000
$ws.Range("I11").Value = 50014000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 50014000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -50013861
$ws.Range("N11").ClearContents()
$ws.Range("H92").Value = 8166
$ws.Range("J92").Value = 8166
$ws.Range("L92").Value = 8166
$ws.Range("N92").Value = -11910
$ws.Range("H99").Value = 15116.667
$ws.Range("J99").Value = 19988.889
$ws.Range("L99").Value = 19988.889
$ws.Range("N99").Value = -24480.889
$ws.Range("H122").Value = 2750.5
$ws.Range("I122").Value = 3170
$ws.Range("J122").Value = 2051.3333
$ws.Range("K122").Value = 9510
$ws.Range("L122").Value = 6153.999899999999
$ws.Range("M122").Value = -7060
$ws.Range("N122").Value = -11053.9999
$ws.Range("H126").Value = 8250.375
$ws.Range("I126").Value = 12050.6
$ws.Range("J126").Value = 1916.6666
$ws.Range("K126").Value = 36151.8
$ws.Range("L126").Value = 5749.9998
$ws.Range("M126").Value = -33681.8
$ws.Range("N126").Value = -10689.9998
$ws.Range("H132").Value = 2622.2068
$ws.Range("I132").Value = 2034.6666
$ws.Range("J132").Value = 3583.6365
$ws.Range("K132").Value = 6103.9998
$ws.Range("L132").Value = 10750.9095
$ws.Range("M132").Value = -3573.9998
$ws.Range("N132").Value = -15810.9095
$ws.Range("H134").Value = 24800
$ws.Range("J134").Value = 24800
$ws.Range("L134").Value = 74400
$ws.Range("N134").Value = -79470
$ws.Range("H135").Value = 32660.637
$ws.Range("J135").Value = 32660.637
$ws.Range("L135").Value = 32660.637
$ws.Range("N135").Value = -42800.637

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 66853.836
$ws.Range("I12").Value = 107
$ws.Range("J12").Value = 80203.2
$ws.Range("K12").Value = 107
$ws.Range("L12").Value = 80203.2
$ws.Range("M12").Value = 63
$ws.Range("N12").Value = -80543.2
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H40").Value = 4519.647
$ws.Range("I40").Value = 4519.647
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4519.647
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4383.647
$ws.Range("N40").ClearContents()
$ws.Range("H136").Value = 2579.3794
$ws.Range("I136").Value = 2078.6086
$ws.Range("J136").Value = 4499
$ws.Range("K136").Value = 6235.825800000001
$ws.Range("L136").Value = 13497
$ws.Range("M136").Value = -3685.825800000001
$ws.Range("N136").Value = -18597

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9524802
$ws.Range("I122").Value = 14286206
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 42858618
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -42856168
$ws.Range("N122").Value = -10885
$ws.Range("H132").Value = 2023.4242
$ws.Range("I132").Value = 1446.421
$ws.Range("J132").Value = 2806.5
$ws.Range("K132").Value = 4339.263
$ws.Range("L132").Value = 8419.5
$ws.Range("M132").Value = -1809.263
$ws.Range("N132").Value = -13479.5
$ws.Range("H136").Value = 21995.818
$ws.Range("I136").Value = 85562.69500000001
$ws.Range("J136").Value = 2320.3572
$ws.Range("K136").Value = 256688.085
$ws.Range("L136").Value = 6961.071599999999
$ws.Range("M136").Value = -254138.085
$ws.Range("N136").Value = -12061.0716

Write-Output "Applied all Masamune_Profits updates."